$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false
$ws = $wb.Worksheets.Item("Book Hoard")

# ---------------------------------------------------------------------------
# 1. Header-row reorder.
#    current_language / original_language move from (P,Q) to just before
#    "topic" (now F,G); fraction_complete moves from X to just before
#    "market_value" (now N). Net effect: columns F..X take on a new set of
#    header labels while A-E and Y onward are untouched.
# ---------------------------------------------------------------------------
$n = 48
$vals = @()
for ($c = 1; $c -le $n; $c++) {
    $vals += $ws.Cells.Item(1, $c).Value2
}

# perm[i] = 0-based source column (in the original layout) that now supplies
# the header text for destination column i (0-based).
$perm = @(0,1,2,3,4,15,16,5,6,7,8,9,10,23,11,12,13,14,17,18,19,20,21,22,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47)

for ($c = 1; $c -le $n; $c++) {
    $ws.Cells.Item(1, $c).Value = $vals[$perm[$c - 1]]
}

# ---------------------------------------------------------------------------
# 2. Header font: bold 11pt -> bold 9pt (affects every header cell, which all
#    share the bold font).
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).Font.Size = 9

# ---------------------------------------------------------------------------
# 3. Column width resizing (post reflow / re-autofit of the header row).
#    ColumnWidth values are chosen so the engine's internal rounding lands as
#    close as possible to the target stored widths.
# ---------------------------------------------------------------------------
$widths = @{
    6  = 16.333333333333332
    7  = 16.666666666666668
    8  = 26.666666666666668
    9  = 25.166666666666668
    10 = 9.333333333333334
    11 = 13.333333333333334
    13 = 17.833333333333332
    14 = 16.5
    15 = 12.0
    16 = 70.66666666666667
    17 = 28.666666666666668
    18 = 22.5
    19 = 68.83333333333333
    20 = 12.333333333333334
    21 = 14.833333333333334
    22 = 16.833333333333332
    23 = 20.5
    24 = 31.166666666666668
    25 = 15.833333333333334
}
foreach ($c in $widths.Keys) {
    $ws.Columns.Item($c).ColumnWidth = $widths[$c]
}

# Extend the explicit-width range from 27-29 to 27-32 (still default width).
$ws.Columns.Item(30).ColumnWidth = 8.333333333333334
$ws.Columns.Item(31).ColumnWidth = 8.333333333333334
$ws.Columns.Item(32).ColumnWidth = 8.333333333333334

# ---------------------------------------------------------------------------
# 4. Remove the empty "Sheet" tab, leaving "Book Hoard" as the only sheet.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet").Delete()

# Re-activate the sole remaining sheet so its tabSelected/activeTab bookkeeping
# is correctly reset now that it is alone in the workbook.
$wb.Worksheets.Item("Book Hoard").Activate()

Write-Output "done"
